# Generate Report for Handback
# Updates the localization-status report: status text, handback timestamps,
# clears the stale "handback file not latest" error once everything is back
# in sync, and widens/narrows the affected columns accordingly.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-31 07:46:34"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-31 07:46:51"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(16).ColumnWidth = 12.8333

Write-Host "Report regenerated for handback."
